$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain-text, dot-grouped numbers (e.g. "43.388.51")
# and also simple decimals (e.g. "302.82"). Excel auto-converts plain decimal-
# looking text into a floating point Number, which would corrupt the intended
# text value (losing trailing zeros / exact formatting). For just those cells,
# pre-set NumberFormat to Text ("@") so the assigned string is kept verbatim.

$ws.Range("D2").Value = "43.388.51"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "2.334.91"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.82"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.19"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E7").Value = "  -0.77%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.504"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.77"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.47"
$ws.Range("E11").Value = "  +7.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0799"
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.92"
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("D15").Value = "2.695.24"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").Value = "2.332.96"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").Value = "43.339.06"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.81"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").Value = "0.0₃0902"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.09"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.06"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.54"
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("E24").Value = "  +4.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.08"
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  +7.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.55"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.14"
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.40"
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.03"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.92"
$ws.Range("E34").Value = "  +5.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.49"
$ws.Range("E35").Value = "  -7.40%  "
$ws.Range("E36").Value = "  +2.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.34"
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.79"
$ws.Range("E40").Value = "  +1.93%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "1.991.69"
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.75"
$ws.Range("E43").Value = "  +6.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0283"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.43"
$ws.Range("E45").Value = "  +4.26%  "
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.90"
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("D49").Value = "2.563.62"
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.11"
$ws.Range("E50").Value = "  +0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.87"
$ws.Range("E51").Value = "  +0.80%  "
